$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting existing rows 30-75 down to 31-76.
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the latest weekly data entry.
$ws.Cells.Item(30, 1).Value = 3
$ws.Cells.Item(30, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(30, 3).Value = "Coquimbo"
$ws.Cells.Item(30, 4).Value = 44792
$ws.Cells.Item(30, 5).Value = 5
$ws.Cells.Item(30, 6).Value = 100112035
$ws.Cells.Item(30, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 50
$ws.Cells.Item(30, 11).Value = 15000
$ws.Cells.Item(30, 12).Value = 15000
$ws.Cells.Item(30, 13).Value = 15000
$ws.Cells.Item(30, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(30, 16).Value = 1000
$ws.Cells.Item(30, 17).Value = 15
$ws.Cells.Item(30, 18).Value = "Hortaliza"
